$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.167.09"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "1.655.30"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.87"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5302"
$ws.Range("E6").Value = "  +2.05%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2624"
$ws.Range("E8").Value = "  +0.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06323"
$ws.Range("E9").Value = "  +1.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.39"
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07807"
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.520"
$ws.Range("E12").Value = "  +1.40%  "
$ws.Range("D13").Value = "1.629.81"
$ws.Range("E13").Value = "  -2.92%  "
$ws.Range("D14").Value = "1.882.90"
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5491"
$ws.Range("E15").Value = "  +1.34%  "
$ws.Range("D16").Value = "0.0₅8164"
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.39"
$ws.Range("E17").Value = "  +1.35%  "
$ws.Range("D18").Value = "26.131.58"
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.600"
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.70"
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("E22").Value = "  +0.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.004"
$ws.Range("E23").Value = "  +0.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.006"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.10"
$ws.Range("E25").Value = "  +4.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1223"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.203"
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.98"
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.472"
$ws.Range("E29").Value = "  +5.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05707"
$ws.Range("E30").Value = "  -3.39%  "
$ws.Range("E31").Value = "  +0.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.549"
$ws.Range("E32").Value = "  +1.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.265"
$ws.Range("E33").Value = "  +1.27%  "
$ws.Range("E34").Value = "  +4.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.802"
$ws.Range("E35").Value = "  +1.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.420"
$ws.Range("E36").Value = "  +0.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9489"
$ws.Range("E37").Value = "  +0.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5727"
$ws.Range("E38").Value = "  +0.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01608"
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.797"
$ws.Range("E40").Value = "  -0.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8498"
$ws.Range("E41").Value = "  +0.91%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "104.44"
$ws.Range("E42").Value = "  +3.92%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.005"
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("D44").Value = "1.039.32"
$ws.Range("E44").Value = "  +4.01%  "
$ws.Range("D45").Value = "1.795.71"
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.76"
$ws.Range("E46").Value = "  +0.58%  "
$ws.Range("E47").Value = "  -2.38%  "
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4352"
$ws.Range("E49").Value = "  +1.41%  "
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.858"
$ws.Range("E51").Value = "  +0.21%  "
